# Generate Report for Handoff
# Replaces the stale GUID-named handoff artifacts with the newly generated
# ones, refreshes the handoff/generate timestamps, and clears out the
# (not-yet-existing) handback target/file/datetime columns for each locale.

$wb = $excel.ActiveWorkbook

$oldGuid = "1958bf87-a8aa-4e62-96af-226375b5b202"
$newGuid = "637f1bc9-7f09-4b97-a488-a75674b08792"

$oldHash = "969ad0aaa490cffa3117e57ca16025d8cdb15281"
$newHash = "c919f550f49e4bfdc699621db79d27461e8a933d"

# 1) Swap every occurrence of the old GUID for the new one across all sheets
#    (source file names, handoff xliff names, hyperlink targets already point
#    to the right branch so only the displayed/cell text needs updating).
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($oldGuid, $newGuid)
    $ws.Cells.Replace($oldHash, $newHash)
}

# 2) Refresh hyperlink display text so it matches the new cell values.
foreach ($ws in $wb.Worksheets) {
    foreach ($h in $ws.Hyperlinks) {
        $h.TextToDisplay = $h.TextToDisplay().Replace($oldGuid, $newGuid)
    }
}

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# 3) Bump the generation / handoff timestamps for the new run.
$wsOverview.Range("G2").Value = "2016-08-25 00:56:33"
$wsZhCn.Range("H2").Value = "2016-08-25 00:56:28"
$wsDeDe.Range("H2").Value = "2016-08-25 00:56:33"

# 4) The handback hasn't happened yet for this run: clear the target file /
#    handback file columns and reset the handback datetime to the "unset"
#    sentinel value for both locales.
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$I$2') {
            $h.Delete()
        }
    }
    $ws.Range("I2").Value = ""
    $ws.Range("J2").Value = ""
    $ws.Range("K2").Value = "0001-01-01 00:00:00"
}

# 5) Narrow the now much shorter "Latest Target File" / "Latest Handback
#    File" columns.
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Range("I1").ColumnWidth = 18.6506053379604
    $ws.Range("J1").ColumnWidth = 21.7054770333426
}
